# "Generate Report for handoff" - record a fresh handoff timestamp for the
# 2404b2f5-... source file in both the zh-cn and de-de status sheets.
$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-02-19 05:27:01"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-02-19 05:27:10"
